# Auto-generated Excel COM-interop script
# Applies scheduled-runner price/profit updates to Halicarnassus_Profits sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(2, 8).Value = 4500
$ws.Cells.Item(2, 10).Value = 4500
$ws.Cells.Item(2, 12).Value = 4500
$ws.Cells.Item(2, 14).Value = -4726
$ws.Cells.Item(17, 8).Value = 1761.0714
$ws.Cells.Item(17, 10).Value = 1761.0714
$ws.Cells.Item(17, 12).Value = 5283.2142
$ws.Cells.Item(17, 14).Value = -5619.2142
$ws.Cells.Item(18, 8).Value = 1420
$ws.Cells.Item(18, 10).Value = 0
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 14).Value = ""
$ws.Cells.Item(32, 8).Value = 0
$ws.Cells.Item(32, 10).Value = 0
$ws.Cells.Item(32, 12).Value = 0
$ws.Cells.Item(32, 14).Value = ""
$ws.Cells.Item(40, 8).Value = 7445.2
$ws.Cells.Item(40, 9).Value = 5864.6665
$ws.Cells.Item(40, 11).Value = 5864.6665
$ws.Cells.Item(40, 13).Value = -5689.6665
$ws.Cells.Item(43, 8).Value = 3500
$ws.Cells.Item(43, 10).Value = 3500
$ws.Cells.Item(43, 12).Value = 3500
$ws.Cells.Item(43, 14).Value = -3638
$ws.Cells.Item(99, 8).Value = 4285.5
$ws.Cells.Item(99, 9).Value = 1714.5
$ws.Cells.Item(99, 11).Value = 5143.5
$ws.Cells.Item(99, 13).Value = -3645.5
$ws.Cells.Item(113, 8).Value = 5174
$ws.Cells.Item(113, 9).Value = 2875.2
$ws.Cells.Item(113, 11).Value = 2875.2
$ws.Cells.Item(113, 13).Value = 378.8000000000002
$ws.Cells.Item(135, 8).Value = 877.5
$ws.Cells.Item(135, 10).Value = 863
$ws.Cells.Item(135, 12).Value = 7767
$ws.Cells.Item(135, 14).Value = -12837
$ws.Cells.Item(137, 8).Value = 1798.2307
$ws.Cells.Item(137, 9).Value = 1064.2222
$ws.Cells.Item(137, 10).Value = 3449.75
$ws.Cells.Item(137, 11).Value = 3192.6666
$ws.Cells.Item(137, 12).Value = 10349.25
$ws.Cells.Item(137, 13).Value = -642.6665999999996
$ws.Cells.Item(137, 14).Value = -15449.25
$ws.Cells.Item(138, 8).Value = 3039.1667
$ws.Cells.Item(138, 9).Value = 990.375
$ws.Cells.Item(138, 10).Value = 4678.2
$ws.Cells.Item(138, 11).Value = 2971.125
$ws.Cells.Item(138, 12).Value = 14034.6
$ws.Cells.Item(138, 13).Value = 2168.875
$ws.Cells.Item(138, 14).Value = -24314.6
$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(2, 8).Value = 976.5833
$ws.Cells.Item(2, 9).Value = 980.4545000000001
$ws.Cells.Item(2, 11).Value = 980.4545000000001
$ws.Cells.Item(2, 13).Value = -867.4545000000001
$ws.Cells.Item(31, 8).Value = 17864.75
$ws.Cells.Item(31, 9).Value = 17864.75
$ws.Cells.Item(31, 11).Value = 17864.75
$ws.Cells.Item(31, 13).Value = -17570.75
$ws.Cells.Item(61, 8).Value = 3047.2666
$ws.Cells.Item(61, 9).Value = 2622.1072
$ws.Cells.Item(61, 11).Value = 2622.1072
$ws.Cells.Item(61, 13).Value = -2410.1072
$ws.Cells.Item(116, 8).Value = 976.5833
$ws.Cells.Item(116, 9).Value = 980.4545000000001
$ws.Cells.Item(116, 11).Value = 980.4545000000001
$ws.Cells.Item(116, 13).Value = 1313.5455
$ws.Cells.Item(136, 8).Value = 3047.2666
$ws.Cells.Item(136, 9).Value = 2622.1072
$ws.Cells.Item(136, 11).Value = 7866.321599999999
$ws.Cells.Item(136, 13).Value = -5316.321599999999
$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(3, 8).Value = 976.5833
$ws.Cells.Item(3, 9).Value = 980.4545000000001
$ws.Cells.Item(3, 11).Value = 980.4545000000001
$ws.Cells.Item(3, 13).Value = -866.4545000000001
$ws.Cells.Item(46, 8).Value = 20000
$ws.Cells.Item(46, 10).Value = 20000
$ws.Cells.Item(46, 12).Value = 20000
$ws.Cells.Item(46, 14).Value = -20596
$ws.Cells.Item(99, 8).Value = 1590.409
$ws.Cells.Item(99, 9).Value = 1341.9412
$ws.Cells.Item(99, 11).Value = 1341.9412
$ws.Cells.Item(99, 13).Value = 156.0588
$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(22, 8).Value = 4500
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).Value = ""
$ws.Cells.Item(74, 8).Value = 41608.285
$ws.Cells.Item(74, 10).Value = 41608.285
$ws.Cells.Item(74, 12).Value = 41608.285
$ws.Cells.Item(74, 14).Value = -43356.285
$ws.Cells.Item(77, 8).Value = 41608.285
$ws.Cells.Item(77, 10).Value = 41608.285
$ws.Cells.Item(77, 12).Value = 124824.855
$ws.Cells.Item(77, 14).Value = -133560.855
$ws.Cells.Item(99, 8).Value = 2226.375
$ws.Cells.Item(99, 9).Value = 1012
$ws.Cells.Item(99, 11).Value = 1012
$ws.Cells.Item(99, 13).Value = 486
$ws.Cells.Item(126, 8).Value = 2226.375
$ws.Cells.Item(126, 9).Value = 1012
$ws.Cells.Item(126, 11).Value = 3036
$ws.Cells.Item(126, 13).Value = -566
$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(22, 8).Value = 260
$ws.Cells.Item(22, 10).Value = 295
$ws.Cells.Item(22, 12).Value = 885
$ws.Cells.Item(22, 14).Value = -1223
$ws.Cells.Item(27, 8).Value = 260
$ws.Cells.Item(27, 10).Value = 295
$ws.Cells.Item(27, 12).Value = 885
$ws.Cells.Item(27, 14).Value = -1089
$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(11, 8).Value = 12186710
$ws.Cells.Item(11, 10).Value = 10001300
$ws.Cells.Item(11, 12).Value = 10001300
$ws.Cells.Item(11, 14).Value = -10001578
$ws.Cells.Item(102, 8).Value = 2453.6191
$ws.Cells.Item(102, 9).Value = 2326.3
$ws.Cells.Item(102, 11).Value = 2326.3
$ws.Cells.Item(102, 13).Value = -704.3000000000002
$ws.Cells.Item(113, 8).Value = 5618.1665
$ws.Cells.Item(113, 9).Value = 3427.5
$ws.Cells.Item(113, 10).Value = 9999.5
$ws.Cells.Item(113, 11).Value = 3427.5
$ws.Cells.Item(113, 12).Value = 9999.5
$ws.Cells.Item(113, 13).Value = -1257.5
$ws.Cells.Item(113, 14).Value = -14339.5
$ws.Cells.Item(126, 8).Value = 1762.4
$ws.Cells.Item(126, 9).Value = 1953
$ws.Cells.Item(126, 11).Value = 5859
$ws.Cells.Item(126, 13).Value = -3389
$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(55, 8).Value = 1727.4286
$ws.Cells.Item(55, 9).Value = 1818.4
$ws.Cells.Item(55, 11).Value = 1818.4
$ws.Cells.Item(55, 13).Value = -1645.4
$ws.Cells.Item(68, 8).Value = 8750
$ws.Cells.Item(68, 9).Value = 0
$ws.Cells.Item(68, 10).Value = 8750
$ws.Cells.Item(68, 11).Value = 0
$ws.Cells.Item(68, 12).Value = 8750
$ws.Cells.Item(68, 13).Value = ""
$ws.Cells.Item(68, 14).Value = -10248
$ws.Cells.Item(71, 8).Value = 8750
$ws.Cells.Item(71, 9).Value = 0
$ws.Cells.Item(71, 10).Value = 8750
$ws.Cells.Item(71, 11).Value = 0
$ws.Cells.Item(71, 12).Value = 43750
$ws.Cells.Item(71, 13).Value = ""
$ws.Cells.Item(71, 14).Value = -51238
$ws.Cells.Item(82, 8).Value = 3800.1428
$ws.Cells.Item(82, 10).Value = 6000
$ws.Cells.Item(82, 12).Value = 6000
$ws.Cells.Item(82, 14).Value = -6722
$ws.Cells.Item(85, 8).Value = 3800.1428
$ws.Cells.Item(85, 10).Value = 6000
$ws.Cells.Item(85, 12).Value = 6000
$ws.Cells.Item(85, 14).Value = -8496
$ws.Cells.Item(136, 8).Value = 3000
$ws.Cells.Item(136, 9).Value = 3000
$ws.Cells.Item(136, 11).Value = 9000
$ws.Cells.Item(136, 13).Value = -6450
$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(81, 8).Value = 3999.8
$ws.Cells.Item(81, 9).Value = 1666.3334
$ws.Cells.Item(81, 11).Value = 3332.6668
$ws.Cells.Item(81, 13).Value = -2271.6668
$ws.Cells.Item(84, 8).Value = 3999.8
$ws.Cells.Item(84, 9).Value = 1666.3334
$ws.Cells.Item(84, 11).Value = 16663.334
$ws.Cells.Item(84, 13).Value = -11359.334
$ws.Cells.Item(126, 8).Value = 5466.2104
$ws.Cells.Item(126, 9).Value = 3384.2222
$ws.Cells.Item(126, 11).Value = 10152.6666
$ws.Cells.Item(126, 13).Value = -7682.6666
$ws.Cells.Item(132, 8).Value = 2090.8
$ws.Cells.Item(132, 9).Value = 1863.625
$ws.Cells.Item(132, 11).Value = 5590.875
$ws.Cells.Item(132, 13).Value = -3060.875

Write-Host "Applied all updates"